$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (Submission #8) ---
$ws.Range("B13").Value = 0.97889599999999999
$ws.Range("C13").Value = 42234
$ws.Range("C13").NumberFormat = "m/d/yy"
$ws.Range("D13").Value = "XGB"
$ws.Range("E13").Value = "eta=.03"

# --- Row 14 (Submission #9) ---
$ws.Range("B14").Value = 0.97886200000000001
$ws.Range("C14").Value = 42234
$ws.Range("C14").NumberFormat = "m/d/yy"
$ws.Range("D14").Value = "XGB"
$ws.Range("E14").Value = "eta=.03"

# --- Row 15 (Submission #10) ---
$ws.Range("B15").Value = 0.98378100000000002
$ws.Range("C15").Value = 42238
$ws.Range("C15").NumberFormat = "m/d/yy"
$ws.Range("D15").Value = "XGB"
$ws.Range("E15").Value = "eta=.03"

# --- Row 16 (Submission #11) ---
$ws.Range("B16").Value = 0.98486200000000002
$ws.Range("C16").Value = 42238
$ws.Range("C16").NumberFormat = "m/d/yy"
$ws.Range("D16").Value = "XGB"
$ws.Range("E16").Value = "eta=.03"

# Preprocessing notes column (H). Enter H14's (longer) text before H13's
# (shorter) text so the shared-string table is built in the same order
# as the source workbook.
$ws.Range("H14").Value = "added date.x.num, date.y.num, date.x.month, date.y.month`ntook out char_1 - char_9"
$ws.Range("H13").Value = "added date.x.num, date.y.num, date.x.month, date.y.month"
$ws.Range("H15").Value = "added group_month feature"
$ws.Range("H16").Value = "added group_month_year feature, removed group_month"

# H13/H14 wrap their long notes and grow taller; H15/H16 stay single-line.
$ws.Range("H13").WrapText = $true
$ws.Range("H14").WrapText = $true
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 45

# Widen the new Preprocessing Notes column to fit its content.
$ws.Columns.Item(8).ColumnWidth = 52.8

# Recalculate so the summary formulas in B1/B2 pick up the new best score.
$wb.Application.Calculate()

# Final selection left on B17, matching where the user's cursor ended up.
$ws.Range("B17").Select() | Out-Null
